$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BMO")

$nl = [char]10
$newResponse = "HTTP/1.1 200 " + $nl + "Content-Type: application/java; charset=windows-1252" + $nl + "Content-Length: 13" + $nl + "Date: Thu, 13 Jun 2019 06:14:34 GMT" + $nl + $nl + "CI/CD Demo..."

# Row 2: Response / PASS / 400
$ws.Range("A2").Value = $newResponse
$ws.Range("A2").WrapText = $true
$ws.Range("B2").Value = "PASS"
$ws.Range("C2").Value = 400

# Row 3: Response / FAIL / 200
$ws.Range("A3").Value = $newResponse
$ws.Range("B3").Value = "FAIL"
$ws.Range("C3").Value = 200

# Match the row height the real Excel autofit produced for the longer wrapped text
$ws.Rows.Item(2).RowHeight = 195
# Row 3 keeps its default (unwrapped) height
$ws.Rows.Item(3).AutoFit()

# Final selection reported in the saved file
$ws.Range("C2").Select() | Out-Null
